$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 63, pushing the existing rows 63-77 down to 64-78.
$ws.Rows.Item(63).Insert()

# Populate the newly inserted row 63 with the new weekly price record.
$ws.Cells.Item(63, 1).Value = 5
$ws.Cells.Item(63, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(63, 3).Value = "Maule"
$ws.Cells.Item(63, 4).Value = 44551
$ws.Cells.Item(63, 5).Value = 7
$ws.Cells.Item(63, 6).Value = 100112030
$ws.Cells.Item(63, 7).Value = "Poroto granado"
$ws.Cells.Item(63, 8).Value = "Sin especificar"
$ws.Cells.Item(63, 9).Value = "Primera"
$ws.Cells.Item(63, 10).Value = 200
$ws.Cells.Item(63, 11).Value = 35000
$ws.Cells.Item(63, 12).Value = 35000
$ws.Cells.Item(63, 13).Value = 35000
$ws.Cells.Item(63, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(63, 15).Value = "Región del Maule"
$ws.Cells.Item(63, 16).Value = 1400
$ws.Cells.Item(63, 17).Value = 25
$ws.Cells.Item(63, 18).Value = "Hortaliza"
